$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to machine-readable codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the lowercase Spanish connector words (de/del/la/las/los/el/y)
# in state/municipality names throughout the data rows
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B9").Value = "Playas De Rosarito"
$ws.Range("B23").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Benemérito De Las Américas"
$ws.Range("B30").Value = "Comitán De Domínguez"
$ws.Range("B68").Value = "Hidalgo Del Parral"
$ws.Range("A71").Value = "Ciudad De México"
$ws.Range("B73").Value = "Cuajimalpa De Morelos"
$ws.Range("A87").Value = "Coahuila De Zaragoza"
$ws.Range("B98").Value = "Villa De Álvarez"
$ws.Range("B109").Value = "San Juan Del Río"
$ws.Range("A115").Value = "Estado De México"
$ws.Range("B115").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B119").Value = "Chapa De Mota"
$ws.Range("B124").Value = "Ecatepec De Morelos"
$ws.Range("B129").Value = "Ixtapan De La Sal"
$ws.Range("B136").Value = "Naucalpan De Juárez"
$ws.Range("B142").Value = "San Felipe Del Progreso"
$ws.Range("B143").Value = "San Martín De Las Pirámides"
$ws.Range("B151").Value = "Tenango Del Valle"
$ws.Range("B157").Value = "Tlalnepantla De Baz"
$ws.Range("B161").Value = "Valle De Bravo"
$ws.Range("B162").Value = "Valle De Chalco Solidaridad"
$ws.Range("B175").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B178").Value = "Jaral Del Progreso"
$ws.Range("B186").Value = "San Francisco Del Rincón"
$ws.Range("B188").Value = "San Luis De La Paz"
$ws.Range("B189").Value = "Silao De La Victoria"
$ws.Range("B192").Value = "Valle De Santiago"
$ws.Range("B196").Value = "Acapulco De Juárez"
$ws.Range("B198").Value = "Ajuchitlán Del Progreso"
$ws.Range("B199").Value = "Alcozauca De Guerrero"
$ws.Range("B204").Value = "Atoyac De Álvarez"
$ws.Range("B205").Value = "Ayutla De Los Libres"
$ws.Range("B207").Value = "Buenavista De Cuéllar"
$ws.Range("B208").Value = "Chilapa De Álvarez"
$ws.Range("B209").Value = "Chilpancingo De Los Bravo"
$ws.Range("B210").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B211").Value = "Cochoapa El Grande"
$ws.Range("B215").Value = "Coyuca De Benítez"
$ws.Range("B216").Value = "Coyuca De Catalán"
$ws.Range("B219").Value = "Cuetzala Del Progreso"
$ws.Range("B224").Value = "Iguala De La Independencia"
$ws.Range("B226").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B227").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B230").Value = "Mártir De Cuilapan"
$ws.Range("B240").Value = "Taxco De Alarcón"
$ws.Range("B243").Value = "Tepecoacuilco De Trujano"
$ws.Range("B248").Value = "Tlapa De Comonfort"
$ws.Range("B250").Value = "Técpan De Galeana"
$ws.Range("B253").Value = "Zihuatanejo De Azueta"
$ws.Range("B260").Value = "Agua Blanca De Iturbide"
$ws.Range("B261").Value = "Atotonilco El Grande"
$ws.Range("B263").Value = "Cuautepec De Hinojosa"
$ws.Range("B273").Value = "Pachuca De Soto"
$ws.Range("B279").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B280").Value = "Tezontepec De Aldama"
$ws.Range("B283").Value = "Tula De Allende"
$ws.Range("B284").Value = "Tulancingo De Bravo"
$ws.Range("B286").Value = "Zacualtipán De Ángeles"
$ws.Range("B287").Value = "Zapotlán De Juárez"
$ws.Range("B292").Value = "Acatlán De Juárez"
$ws.Range("B295").Value = "Atemajac De Brizuela"
$ws.Range("B296").Value = "Atotonilco El Alto"
$ws.Range("B298").Value = "Autlán De Navarro"
$ws.Range("B302").Value = "Cañadas De Obregón"
$ws.Range("B308").Value = "Cuautitlán De García Barragán"
$ws.Range("B314").Value = "Encarnación De Díaz"
$ws.Range("B317").Value = "Huejuquilla El Alto"
$ws.Range("B318").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B319").Value = "Ixtlahuacán Del Río"
$ws.Range("B321").Value = "Jilotlán De Los Dolores"
$ws.Range("B326").Value = "Lagos De Moreno"
$ws.Range("B329").Value = "Ojuelos De Jalisco"
$ws.Range("B334").Value = "San Juan De Los Lagos"
$ws.Range("B335").Value = "San Juanito De Escobedo"
$ws.Range("B337").Value = "San Martín De Bolaños"
$ws.Range("B338").Value = "San Miguel El Alto"
$ws.Range("B339").Value = "San Sebastián Del Oeste"
$ws.Range("B341").Value = "Talpa De Allende"
$ws.Range("B342").Value = "Tamazula De Gordiano"
$ws.Range("B347").Value = "Teocuitatlán De Corona"
$ws.Range("B348").Value = "Tepatitlán De Morelos"
$ws.Range("B349").Value = "Tizapán El Alto"
$ws.Range("B350").Value = "Tlajomulco De Zúñiga"
$ws.Range("B355").Value = "Unión De Tula"
$ws.Range("B358").Value = "Zacoalco De Torres"
$ws.Range("B360").Value = "Zapotitlán De Vadillo"
$ws.Range("B361").Value = "Zapotlán El Grande"
$ws.Range("A363").Value = "Michoacán De Ocampo"
$ws.Range("B377").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B379").Value = "Cojumatlán De Régules"
$ws.Range("B429").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B458").Value = "Puente De Ixtla"
$ws.Range("B469").Value = "Bahía De Banderas"
$ws.Range("B473").Value = "Ixtlán Del Río"
$ws.Range("B480").Value = "Santa María Del Oro"
$ws.Range("B492").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B498").Value = "Chalcatongo De Hidalgo"
$ws.Range("B500").Value = "Coicoyán De Las Flores"
$ws.Range("B501").Value = "Constancia Del Rosario"
$ws.Range("B502").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B503").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B504").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B505").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B506").Value = "Huajuapan De León"
$ws.Range("B507").Value = "Ixtlán De Juárez"
$ws.Range("B511").Value = "Mariscala De Juárez"
$ws.Range("B513").Value = "Mazatlán Villa De Flores"
$ws.Range("B514").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B515").Value = "Nejapa De Madero"
$ws.Range("B516").Value = "Oaxaca De Juárez"
$ws.Range("B517").Value = "Ocotlán De Morelos"
$ws.Range("B518").Value = "Putla Villa De Guerrero"
$ws.Range("B540").Value = "San José Del Peñasco"
$ws.Range("B543").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B556").Value = "San Juan Del Estado"
$ws.Range("B573").Value = "San Miguel Del Puerto"
$ws.Range("B589").Value = "San Pedro El Alto"
$ws.Range("B601").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B606").Value = "Santa Inés Del Monte"
$ws.Range("B646").Value = "Santo Domingo De Morelos"
$ws.Range("B650").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B651").Value = "Tlacolula De Matamoros"
$ws.Range("B652").Value = "Totontepec Villa De Morelos"
$ws.Range("B653").Value = "Villa Sola De Vega"
$ws.Range("B654").Value = "Villa De Etla"
$ws.Range("B655").Value = "Villa De Zaachila"
$ws.Range("B658").Value = "Zimatlán De Álvarez"
$ws.Range("B684").Value = "Huehuetlán El Grande"
$ws.Range("B688").Value = "Izúcar De Matamoros"
$ws.Range("B693").Value = "Los Reyes De Juárez"
$ws.Range("B694").Value = "Mazapiltepec De Juárez"
$ws.Range("B698").Value = "Palmar De Bravo"
$ws.Range("B713").Value = "Tepanco De López"
$ws.Range("B715").Value = "Tepexi De Rodríguez"
$ws.Range("B716").Value = "Tetela De Ocampo"
$ws.Range("B720").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B736").Value = "Amealco De Bonfil"
$ws.Range("B737").Value = "Cadereyta De Montes"
$ws.Range("B742").Value = "Jalpan De Serra"
$ws.Range("B743").Value = "Pinal De Amoles"
$ws.Range("B745").Value = "San Juan Del Río"
$ws.Range("B753").Value = "Axtla De Terrazas"
$ws.Range("B759").Value = "Villa De Ramos"
$ws.Range("B806").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B807").Value = "San Pablo Del Monte"
$ws.Range("B808").Value = "Tetla De La Solidaridad"
$ws.Range("A816").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B819").Value = "Amatlán De Los Reyes"
$ws.Range("B827").Value = "Boca Del Río"
$ws.Range("B833").Value = "Cosamaloapan De Carpio"
$ws.Range("B843").Value = "Ignacio De La Llave"
$ws.Range("B846").Value = "Ixhuatlán Del Café"
$ws.Range("B852").Value = "Juchique De Ferrer"
$ws.Range("B855").Value = "Lerdo De Tejada"
$ws.Range("B858").Value = "Martínez De La Torre"
$ws.Range("B859").Value = "Medellín De Bravo"
$ws.Range("B861").Value = "Mixtla De Altamirano"
$ws.Range("B868").Value = "Poza Rica De Hidalgo"
$ws.Range("B875").Value = "Soledad De Doblado"
$ws.Range("B893").Value = "Vega De Alatorre"
$ws.Range("B923").Value = "Jiménez Del Teul"
$ws.Range("B927").Value = "Moyahua De Estrada"
$ws.Range("B928").Value = "Nochistlán De Mejía"
$ws.Range("B932").Value = "Teúl De González Ortega"
$ws.Range("B933").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B937").Value = "Villa De Cos"

# Normalize grand-total label casing
$ws.Range("A941").Value = "Total"

# Remove the trailing footnote rows (943-947) that are no longer part of the clean dataset
$ws.Range("A943:D947").ClearContents()
